$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (id 1): add a new hyperlink in column G pointing at the deployed movie-recommender app ---
$ws.Range("G2").Value = "https://movies-u-like.herokuapp.com"
$ws.Hyperlinks.Add($ws.Range("G2"), "https://movies-u-like.herokuapp.com") | Out-Null
$ws.Range("G2").Style = "Hyperlink"

# --- Row 3 (id 2): add a new hyperlink in column G pointing at the portfolio website ---
$ws.Range("G3").Value = "www.navrozlamba.com"
$ws.Hyperlinks.Add($ws.Range("G3"), "www.navrozlamba.com") | Out-Null
$ws.Range("G3").Style = "Hyperlink"

# --- Row 5 (id 4): replace the existing hyperlink/text with the new deployed Airbnb app link ---
$ws.Range("G5").Hyperlinks.Delete()
$ws.Range("G5").Value = "https://airbnb-la.herokuapp.com/"
$ws.Hyperlinks.Add($ws.Range("G5"), "https://airbnb-la.herokuapp.com/") | Out-Null
$ws.Range("G5").Style = "Hyperlink"

# --- Row 7 (id 6): this row now documents the Airbus data-storytelling project (images added) instead
#     of the stray blog-post row; drop its blog-post hyperlink and rewrite title/description/id ---
$ws.Range("G7").Hyperlinks.Delete()
$ws.Range("G7").ClearContents()

$ws.Range("B7").Value = "Airbus, the new King of the Skies?! | Data Storytelling"
$ws.Range("C7").Value = "Built a narrative around a set of data, and accompanying visualizations to help convey the sense of competition between the two biggest rivalries of all times; Boeing and Airbus."
$ws.Range("E7").Value = "project6"

# --- restore the active selection like the authored workbook ---
$ws.Range("G15").Select()

Write-Output "done"
